# This script re-generates the input-file stimuli table by remapping each
# trial row to a different source row (drawn from one of the 20 base designs)
# and shifting trial_total (col F) by +81 so the numbering continues for the
# next duplicated block of 1000-subject input files.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2 = 17
    3 = 32
    4 = 20
    5 = 13
    6 = 11
    7 = 39
    8 = 23
    9 = 8
    10 = 30
    11 = 16
    12 = 26
    13 = 40
    14 = 7
    15 = 24
    16 = 5
    17 = 21
    18 = 35
    19 = 15
    20 = 14
    21 = 25
    22 = 19
    23 = 10
    24 = 9
    25 = 2
    26 = 38
    27 = 6
    28 = 12
    29 = 37
    30 = 22
    31 = 18
    32 = 34
    33 = 4
    34 = 33
    35 = 41
    36 = 29
    37 = 31
    38 = 36
    39 = 27
    40 = 28
    41 = 3
}

$cols = @("G","H","I","K","L","M","N","O","P","Q","R","S","T","U","V")

# Snapshot the "before" values for every relevant cell first, since the
# source and destination row ranges overlap (this is a permutation).
$snapshot = @{}
foreach ($r in 2..41) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $snapshot[$addr] = $ws.Range($addr).Value()
    }
    $snapshot["F$r"] = $ws.Range("F$r").Value()
}

foreach ($r in 2..41) {
    $src = $rowMap[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $snapshot["$c$src"]
    }
    $ws.Range("F$r").Value = $snapshot["F$r"] + 81
}
